# "get projects from login user: worked"
#
# 1) Slide 2 ("schemas"), the "projects schema" textbox (TextBox 5):
#      - add a new bullet paragraph "-url : string //feature url" right
#        after the "-is_shared: false // boolean" line (keeping the same
#        text colour as the rest of that paragraph), and
#      - grow the shape's height to fit the extra line.
# 2) Slide 4, the small "Load/edit " caption textbox (TextBox 26):
#      - widen it slightly.

$EMU_PER_POINT = 12700

$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Slide 2 - "projects schema" textbox
# ---------------------------------------------------------------------
$slide2 = $p.Slides.Item(2)
$projSchema = $slide2.Shapes.Item(3)   # "TextBox 5" - projects schema

$tr = $projSchema.TextFrame.TextRange

# Locate the paragraph that ends in "... // boolean" and append a new
# paragraph after it, reusing the formatting already on that run so the
# new text inherits the same (theme) text colour.
$fullText = $tr.Text
$boolIdx = $fullText.IndexOf("is_shared: false // boolean")
$paraCount = $tr.Paragraphs().Count
for ($i = 1; $i -le $paraCount; $i++) {
    $candidate = $tr.Paragraphs($i, 1)
    if ($candidate.Text -like "*is_shared: false // boolean*") {
        $boolPara = $candidate
    }
}

$newParaText = [char]13 + "-url : string //feature url"
$null = $boolPara.InsertAfter($newParaText)

# Split the freshly inserted paragraph into separate runs matching the
# original authoring (each "word" as its own run) while keeping the
# inherited colour/size - re-asserting the (already correct) font size
# is enough to force the run boundary without altering formatting.
$afterText = $tr.Text
$newStart = $afterText.IndexOf("-url : string //feature url")

$seg1 = $tr.Characters($newStart + 1, 1)   # "-"
$seg2 = $tr.Characters($newStart + 2, 3)   # "url"
$seg3 = $tr.Characters($newStart + 5, 20)  # " : string //feature "
$seg4 = $tr.Characters($newStart + 25, 3)  # "url"

$seg1.Font.Size = 14
$seg2.Font.Size = 14
$seg3.Font.Size = 14
$seg4.Font.Size = 14

# Grow the textbox so the new line fits (matches the authored resize).
$projSchema.Height = 2585323 / $EMU_PER_POINT

# ---------------------------------------------------------------------
# 2) Slide 4 - "Load/edit " caption textbox
# ---------------------------------------------------------------------
$slide4 = $p.Slides.Item(4)
$loadEdit = $slide4.Shapes.Item(13)   # "TextBox 26" - "Load/edit "
$loadEdit.Width = 1026160 / $EMU_PER_POINT
